$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns V (SystemStartDate) and W (SystemEndDate)
$ws.Range("V1").Value = "SystemStartDate"
$ws.Range("W1").Value = "SystemEndDate"
$ws.Range("V2").Value = "System.Today() 10:00:00"
$ws.Range("W2").Value = "System.Today 16:00:00"

# Column widths for the two new columns
$ws.Columns("V").ColumnWidth = 40.666666666666664
$ws.Columns("W").ColumnWidth = 28.666666666666664

# T2 loses its custom number format (reverts to the same plain wrap-text
# style used by U2)
$ws.Range("U2").Copy()
$ws.Range("T2").PasteSpecial(-4122)

# Update the view: scroll so column Q is the left-most visible column and
# select W8
[void]$ws.Range("W8").Select()
